# Fruta / hortaliza, semanal
#
# A new weekly record is inserted as row 5 (Macroferia Regional de Talca -
# Arandano (blue), week of 2021-11-25), which pushes all subsequent data
# rows (previously rows 5-30) down by one, to rows 6-31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 5, shifting rows 5:30 down
# to 6:31 (Excel copies the row-above formatting, so column D keeps its
# date number format).
$ws.Rows(5).Insert()

# Populate the newly inserted row 5 with the new weekly entry.
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(5, 3).Value = "Maule"
$ws.Cells.Item(5, 4).Value = 44525
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100101
$ws.Cells.Item(5, 8).Value = "Berries"
$ws.Cells.Item(5, 9).Value = 100101001
$ws.Cells.Item(5, 10).Value = "Arándano (blue)"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 150
$ws.Cells.Item(5, 14).Value = 4000
$ws.Cells.Item(5, 15).Value = 4000
$ws.Cells.Item(5, 16).Value = 4000
$ws.Cells.Item(5, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(5, 18).Value = "Provincia de Linares"
$ws.Cells.Item(5, 19).Value = 2000
$ws.Cells.Item(5, 20).Value = 2
